$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 51

# Row 3
$ws.Range("B3").Value = "<shift>"

# Row 4
$ws.Range("B4").Value = "<sine>"
$ws.Range("C4").Value = 55

# Row 6
$ws.Range("C6").Value = 55

# Row 7
$ws.Range("B7").Value = "<with>"
$ws.Range("C7").Value = 55

# Row 8
$ws.Range("C8").Value = 50

# Row 9
$ws.Range("B9").Value = "<it>"
$ws.Range("C9").Value = 14
